$p = $ppt.ActivePresentation

# Slide 3 ("Histogram and Boxplot") currently has two pictures: the
# histogram (kept) and the boxplot (moved out to its own new slide).
$s3 = $p.Slides.Item(3)

# Locate the boxplot picture by its alt-text/description so this keeps
# working even if shape ordering ever changes.
$boxplotShape = $null
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $shp = $s3.Shapes.Item($i)
    if ($shp.Name -eq "Picture 4") {
        $boxplotShape = $shp
    }
}
if ($boxplotShape -eq $null) {
    $boxplotShape = $s3.Shapes.Item($s3.Shapes.Count)
}

# Copy it to the clipboard before removing it from slide 3.
$boxplotShape.Copy() | Out-Null

# Add the new slide at the end of the deck, using the same
# "Title and Content" layout as the other chart slides.
$s6 = $p.Slides.Add(6, 2)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Boxplot"

# Now remove the boxplot picture from slide 3.
$boxplotShape.Delete()

# Paste the boxplot picture (same embedded image, same position/size)
# onto the new slide and restore its original shape name.
$s6.Shapes.Paste() | Out-Null
$pastedPic = $s6.Shapes.Item($s6.Shapes.Count)
$pastedPic.Name = "Picture 3"
